# Update column F (dSF) values on Sheet1 to reflect repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    2  = 4
    5  = 5
    6  = -6
    7  = 1
    8  = -1
    9  = 6
    13 = -2
    15 = -5
    16 = 1
    17 = 4
    18 = -1
    19 = 2
    20 = -3
    21 = 6
    22 = 8
    23 = -4
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
